# Update with midnight and simplification of the time reading for en and de
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 1-3 keep their original meaning (am / pm / o'clock).
# Rows 4-6 are repurposed: "quarter past"/"half past"/"quarter to" go away,
# replaced by "It is" / "midnight" / "noon".
$ws.Range("A1").Value = "am"
$ws.Range("A2").Value = "pm"
$ws.Range("A3").Value = "o’clock"
$ws.Range("A4").Value = "It is"
$ws.Range("A5").Value = "midnight"
$ws.Range("A6").Value = "noon"

# Rows 7 onward become a plain numeric minute/hour sequence 1..59
# (replaces the old text rows "It is"/"past"/"to" plus the old 1..30 list).
$row = 7
for ($n = 1; $n -le 59; $n++) {
    $ws.Cells.Item($row, 1).Value = $n
    $row++
}

# The (formerly unused) second column now carries an explicit Text number
# format, matching the new style added to the workbook.
$ws.Range("B1:B65").NumberFormat = "@"

# Selection moves to A7 in the edited workbook.
$null = $ws.Range("A7").Select()
